$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 499 (shifts old rows 499-517 down to 501-519)
$ws.Range("A499:A500").EntireRow.Insert()

# Row 499
$ws.Cells.Item(499, 1).Value = 7
$ws.Cells.Item(499, 2).Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(499, 3).Value = 'Ñuble'
$ws.Cells.Item(499, 4).Value = Get-Date -Year 2023 -Month 1 -Day 13 -Hour 0 -Minute 0 -Second 0
$ws.Cells.Item(499, 5).Value = 16
$ws.Cells.Item(499, 6).Value = 100114014
$ws.Cells.Item(499, 7).Value = 'Betarraga'
$ws.Cells.Item(499, 8).Value = 'Sin especificar'
$ws.Cells.Item(499, 9).Value = 'Primera'
$ws.Cells.Item(499, 10).Value = 300
$ws.Cells.Item(499, 11).Value = 700
$ws.Cells.Item(499, 12).Value = 750
$ws.Cells.Item(499, 13).Value = 725
$ws.Cells.Item(499, 14).Value = '$/paquete 5 unidades'
$ws.Cells.Item(499, 15).Value = 'Provincia de Diguillín'
$ws.Cells.Item(499, 16).Value = 145
$ws.Cells.Item(499, 17).Value = 5
$ws.Cells.Item(499, 18).Value = 'Hortaliza'

# Row 500
$ws.Cells.Item(500, 1).Value = 7
$ws.Cells.Item(500, 2).Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(500, 3).Value = 'Ñuble'
$ws.Cells.Item(500, 4).Value = Get-Date -Year 2023 -Month 1 -Day 13 -Hour 0 -Minute 0 -Second 0
$ws.Cells.Item(500, 5).Value = 16
$ws.Cells.Item(500, 6).Value = 100114014
$ws.Cells.Item(500, 7).Value = 'Betarraga'
$ws.Cells.Item(500, 8).Value = 'Sin especificar'
$ws.Cells.Item(500, 9).Value = 'Segunda'
$ws.Cells.Item(500, 10).Value = 200
$ws.Cells.Item(500, 11).Value = 600
$ws.Cells.Item(500, 12).Value = 600
$ws.Cells.Item(500, 13).Value = 600
$ws.Cells.Item(500, 14).Value = '$/paquete 5 unidades'
$ws.Cells.Item(500, 15).Value = 'Provincia de Diguillín'
$ws.Cells.Item(500, 16).Value = 120
$ws.Cells.Item(500, 17).Value = 5
$ws.Cells.Item(500, 18).Value = 'Hortaliza'

# Row 501
$ws.Cells.Item(501, 1).Value = 7
$ws.Cells.Item(501, 2).Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(501, 3).Value = 'Ñuble'
$ws.Cells.Item(501, 4).Value = Get-Date -Year 2021 -Month 8 -Day 17 -Hour 0 -Minute 0 -Second 0
$ws.Cells.Item(501, 5).Value = 16
$ws.Cells.Item(501, 6).Value = 100114014
$ws.Cells.Item(501, 7).Value = 'Betarraga'
$ws.Cells.Item(501, 8).Value = 'Sin especificar'
$ws.Cells.Item(501, 9).Value = 'Primera'
$ws.Cells.Item(501, 10).Value = 120
$ws.Cells.Item(501, 11).Value = 600
$ws.Cells.Item(501, 12).Value = 650
$ws.Cells.Item(501, 13).Value = 625
$ws.Cells.Item(501, 14).Value = '$/paquete 5 unidades'
$ws.Cells.Item(501, 15).Value = 'Región del Maule'
$ws.Cells.Item(501, 16).Value = 125
$ws.Cells.Item(501, 17).Value = 5
$ws.Cells.Item(501, 18).Value = 'Hortaliza'

# Row 502
$ws.Cells.Item(502, 1).Value = 7
$ws.Cells.Item(502, 2).Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(502, 3).Value = 'Ñuble'
$ws.Cells.Item(502, 4).Value = Get-Date -Year 2021 -Month 8 -Day 10 -Hour 0 -Minute 0 -Second 0
$ws.Cells.Item(502, 5).Value = 16
$ws.Cells.Item(502, 6).Value = 100114014
$ws.Cells.Item(502, 7).Value = 'Betarraga'
$ws.Cells.Item(502, 8).Value = 'Sin especificar'
$ws.Cells.Item(502, 9).Value = 'Primera'
$ws.Cells.Item(502, 10).Value = 160
$ws.Cells.Item(502, 11).Value = 600
$ws.Cells.Item(502, 12).Value = 650
$ws.Cells.Item(502, 13).Value = 625
$ws.Cells.Item(502, 14).Value = '$/paquete 5 unidades'
$ws.Cells.Item(502, 15).Value = 'Provincia de Diguillín'
$ws.Cells.Item(502, 16).Value = 125
$ws.Cells.Item(502, 17).Value = 5
$ws.Cells.Item(502, 18).Value = 'Hortaliza'

# Row 503
$ws.Cells.Item(503, 1).Value = 7
$ws.Cells.Item(503, 2).Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(503, 3).Value = 'Ñuble'
$ws.Cells.Item(503, 4).Value = Get-Date -Year 2021 -Month 8 -Day 10 -Hour 0 -Minute 0 -Second 0
$ws.Cells.Item(503, 5).Value = 16
$ws.Cells.Item(503, 6).Value = 100114014
$ws.Cells.Item(503, 7).Value = 'Betarraga'
$ws.Cells.Item(503, 8).Value = 'Sin especificar'
$ws.Cells.Item(503, 9).Value = 'Segunda'
$ws.Cells.Item(503, 10).Value = 120
$ws.Cells.Item(503, 11).Value = 500
$ws.Cells.Item(503, 12).Value = 550
$ws.Cells.Item(503, 13).Value = 525
$ws.Cells.Item(503, 14).Value = '$/paquete 5 unidades'
$ws.Cells.Item(503, 15).Value = 'Provincia de Diguillín'
$ws.Cells.Item(503, 16).Value = 105
$ws.Cells.Item(503, 17).Value = 5
$ws.Cells.Item(503, 18).Value = 'Hortaliza'

# Row 504
$ws.Cells.Item(504, 1).Value = 7
$ws.Cells.Item(504, 2).Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(504, 3).Value = 'Ñuble'
$ws.Cells.Item(504, 4).Value = Get-Date -Year 2021 -Month 7 -Day 30 -Hour 0 -Minute 0 -Second 0
$ws.Cells.Item(504, 5).Value = 16
$ws.Cells.Item(504, 6).Value = 100114014
$ws.Cells.Item(504, 7).Value = 'Betarraga'
$ws.Cells.Item(504, 8).Value = 'Sin especificar'
$ws.Cells.Item(504, 9).Value = 'Primera'
$ws.Cells.Item(504, 10).Value = 120
$ws.Cells.Item(504, 11).Value = 600
$ws.Cells.Item(504, 12).Value = 650
$ws.Cells.Item(504, 13).Value = 625
$ws.Cells.Item(504, 14).Value = '$/paquete 5 unidades'
$ws.Cells.Item(504, 15).Value = 'Provincia de Diguillín'
$ws.Cells.Item(504, 16).Value = 125
$ws.Cells.Item(504, 17).Value = 5
$ws.Cells.Item(504, 18).Value = 'Hortaliza'

# Row 505
$ws.Cells.Item(505, 1).Value = 7
$ws.Cells.Item(505, 2).Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(505, 3).Value = 'Ñuble'
$ws.Cells.Item(505, 4).Value = Get-Date -Year 2021 -Month 7 -Day 30 -Hour 0 -Minute 0 -Second 0
$ws.Cells.Item(505, 5).Value = 16
$ws.Cells.Item(505, 6).Value = 100114014
$ws.Cells.Item(505, 7).Value = 'Betarraga'
$ws.Cells.Item(505, 8).Value = 'Sin especificar'
$ws.Cells.Item(505, 9).Value = 'Segunda'
$ws.Cells.Item(505, 10).Value = 120
$ws.Cells.Item(505, 11).Value = 500
$ws.Cells.Item(505, 12).Value = 550
$ws.Cells.Item(505, 13).Value = 525
$ws.Cells.Item(505, 14).Value = '$/paquete 5 unidades'
$ws.Cells.Item(505, 15).Value = 'Provincia de Diguillín'
$ws.Cells.Item(505, 16).Value = 105
$ws.Cells.Item(505, 17).Value = 5
$ws.Cells.Item(505, 18).Value = 'Hortaliza'

# Row 506
$ws.Cells.Item(506, 1).Value = 7
$ws.Cells.Item(506, 2).Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(506, 3).Value = 'Ñuble'
$ws.Cells.Item(506, 4).Value = Get-Date -Year 2021 -Month 9 -Day 28 -Hour 0 -Minute 0 -Second 0
$ws.Cells.Item(506, 5).Value = 16
$ws.Cells.Item(506, 6).Value = 100114014
$ws.Cells.Item(506, 7).Value = 'Betarraga'
$ws.Cells.Item(506, 8).Value = 'Sin especificar'
$ws.Cells.Item(506, 9).Value = 'Primera'
$ws.Cells.Item(506, 10).Value = 300
$ws.Cells.Item(506, 11).Value = 750
$ws.Cells.Item(506, 12).Value = 800
$ws.Cells.Item(506, 13).Value = 775
$ws.Cells.Item(506, 14).Value = '$/paquete 5 unidades'
$ws.Cells.Item(506, 15).Value = 'Provincia de Diguillín'
$ws.Cells.Item(506, 16).Value = 155
$ws.Cells.Item(506, 17).Value = 5
$ws.Cells.Item(506, 18).Value = 'Hortaliza'

# Row 507
$ws.Cells.Item(507, 1).Value = 7
$ws.Cells.Item(507, 2).Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(507, 3).Value = 'Ñuble'
$ws.Cells.Item(507, 4).Value = Get-Date -Year 2022 -Month 9 -Day 23 -Hour 0 -Minute 0 -Second 0
$ws.Cells.Item(507, 5).Value = 16
$ws.Cells.Item(507, 6).Value = 100114014
$ws.Cells.Item(507, 7).Value = 'Betarraga'
$ws.Cells.Item(507, 8).Value = 'Sin especificar'
$ws.Cells.Item(507, 9).Value = 'Primera'
$ws.Cells.Item(507, 10).Value = 160
$ws.Cells.Item(507, 11).Value = 750
$ws.Cells.Item(507, 12).Value = 800
$ws.Cells.Item(507, 13).Value = 775
$ws.Cells.Item(507, 14).Value = '$/paquete 5 unidades'
$ws.Cells.Item(507, 15).Value = 'Provincia de Diguillín'
$ws.Cells.Item(507, 16).Value = 155
$ws.Cells.Item(507, 17).Value = 5
$ws.Cells.Item(507, 18).Value = 'Hortaliza'

# Row 508
$ws.Cells.Item(508, 1).Value = 7
$ws.Cells.Item(508, 2).Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(508, 3).Value = 'Ñuble'
$ws.Cells.Item(508, 4).Value = Get-Date -Year 2022 -Month 5 -Day 9 -Hour 0 -Minute 0 -Second 0
$ws.Cells.Item(508, 5).Value = 16
$ws.Cells.Item(508, 6).Value = 100114014
$ws.Cells.Item(508, 7).Value = 'Betarraga'
$ws.Cells.Item(508, 8).Value = 'Sin especificar'
$ws.Cells.Item(508, 9).Value = 'Primera'
$ws.Cells.Item(508, 10).Value = 240
$ws.Cells.Item(508, 11).Value = 700
$ws.Cells.Item(508, 12).Value = 750
$ws.Cells.Item(508, 13).Value = 725
$ws.Cells.Item(508, 14).Value = '$/paquete 5 unidades'
$ws.Cells.Item(508, 15).Value = 'Provincia de Diguillín'
$ws.Cells.Item(508, 16).Value = 145
$ws.Cells.Item(508, 17).Value = 5
$ws.Cells.Item(508, 18).Value = 'Hortaliza'

# Row 509
$ws.Cells.Item(509, 1).Value = 7
$ws.Cells.Item(509, 2).Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(509, 3).Value = 'Ñuble'
$ws.Cells.Item(509, 4).Value = Get-Date -Year 2022 -Month 5 -Day 9 -Hour 0 -Minute 0 -Second 0
$ws.Cells.Item(509, 5).Value = 16
$ws.Cells.Item(509, 6).Value = 100114014
$ws.Cells.Item(509, 7).Value = 'Betarraga'
$ws.Cells.Item(509, 8).Value = 'Sin especificar'
$ws.Cells.Item(509, 9).Value = 'Segunda'
$ws.Cells.Item(509, 10).Value = 150
$ws.Cells.Item(509, 11).Value = 600
$ws.Cells.Item(509, 12).Value = 600
$ws.Cells.Item(509, 13).Value = 600
$ws.Cells.Item(509, 14).Value = '$/paquete 5 unidades'
$ws.Cells.Item(509, 15).Value = 'Provincia de Diguillín'
$ws.Cells.Item(509, 16).Value = 120
$ws.Cells.Item(509, 17).Value = 5
$ws.Cells.Item(509, 18).Value = 'Hortaliza'

# Row 510
$ws.Cells.Item(510, 1).Value = 7
$ws.Cells.Item(510, 2).Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(510, 3).Value = 'Ñuble'
$ws.Cells.Item(510, 4).Value = Get-Date -Year 2021 -Month 10 -Day 19 -Hour 0 -Minute 0 -Second 0
$ws.Cells.Item(510, 5).Value = 16
$ws.Cells.Item(510, 6).Value = 100114014
$ws.Cells.Item(510, 7).Value = 'Betarraga'
$ws.Cells.Item(510, 8).Value = 'Sin especificar'
$ws.Cells.Item(510, 9).Value = 'Primera'
$ws.Cells.Item(510, 10).Value = 200
$ws.Cells.Item(510, 11).Value = 750
$ws.Cells.Item(510, 12).Value = 800
$ws.Cells.Item(510, 13).Value = 775
$ws.Cells.Item(510, 14).Value = '$/paquete 5 unidades'
$ws.Cells.Item(510, 15).Value = 'Región del Maule'
$ws.Cells.Item(510, 16).Value = 155
$ws.Cells.Item(510, 17).Value = 5
$ws.Cells.Item(510, 18).Value = 'Hortaliza'

# Row 511
$ws.Cells.Item(511, 1).Value = 7
$ws.Cells.Item(511, 2).Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(511, 3).Value = 'Ñuble'
$ws.Cells.Item(511, 4).Value = Get-Date -Year 2021 -Month 6 -Day 23 -Hour 0 -Minute 0 -Second 0
$ws.Cells.Item(511, 5).Value = 16
$ws.Cells.Item(511, 6).Value = 100114014
$ws.Cells.Item(511, 7).Value = 'Betarraga'
$ws.Cells.Item(511, 8).Value = 'Sin especificar'
$ws.Cells.Item(511, 9).Value = 'Primera'
$ws.Cells.Item(511, 10).Value = 120
$ws.Cells.Item(511, 11).Value = 650
$ws.Cells.Item(511, 12).Value = 700
$ws.Cells.Item(511, 13).Value = 675
$ws.Cells.Item(511, 14).Value = '$/paquete 5 unidades'
$ws.Cells.Item(511, 15).Value = 'Provincia de Diguillín'
$ws.Cells.Item(511, 16).Value = 135
$ws.Cells.Item(511, 17).Value = 5
$ws.Cells.Item(511, 18).Value = 'Hortaliza'

# Row 512
$ws.Cells.Item(512, 1).Value = 7
$ws.Cells.Item(512, 2).Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(512, 3).Value = 'Ñuble'
$ws.Cells.Item(512, 4).Value = Get-Date -Year 2021 -Month 5 -Day 5 -Hour 0 -Minute 0 -Second 0
$ws.Cells.Item(512, 5).Value = 16
$ws.Cells.Item(512, 6).Value = 100114014
$ws.Cells.Item(512, 7).Value = 'Betarraga'
$ws.Cells.Item(512, 8).Value = 'Sin especificar'
$ws.Cells.Item(512, 9).Value = 'Primera'
$ws.Cells.Item(512, 10).Value = 120
$ws.Cells.Item(512, 11).Value = 600
$ws.Cells.Item(512, 12).Value = 650
$ws.Cells.Item(512, 13).Value = 625
$ws.Cells.Item(512, 14).Value = '$/paquete 5 unidades'
$ws.Cells.Item(512, 15).Value = 'Provincia de Diguillín'
$ws.Cells.Item(512, 16).Value = 125
$ws.Cells.Item(512, 17).Value = 5
$ws.Cells.Item(512, 18).Value = 'Hortaliza'

# Row 513
$ws.Cells.Item(513, 1).Value = 7
$ws.Cells.Item(513, 2).Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(513, 3).Value = 'Ñuble'
$ws.Cells.Item(513, 4).Value = Get-Date -Year 2021 -Month 5 -Day 5 -Hour 0 -Minute 0 -Second 0
$ws.Cells.Item(513, 5).Value = 16
$ws.Cells.Item(513, 6).Value = 100114014
$ws.Cells.Item(513, 7).Value = 'Betarraga'
$ws.Cells.Item(513, 8).Value = 'Sin especificar'
$ws.Cells.Item(513, 9).Value = 'Segunda'
$ws.Cells.Item(513, 10).Value = 40
$ws.Cells.Item(513, 11).Value = 500
$ws.Cells.Item(513, 12).Value = 500
$ws.Cells.Item(513, 13).Value = 500
$ws.Cells.Item(513, 14).Value = '$/paquete 5 unidades'
$ws.Cells.Item(513, 15).Value = 'Provincia de Diguillín'
$ws.Cells.Item(513, 16).Value = 100
$ws.Cells.Item(513, 17).Value = 5
$ws.Cells.Item(513, 18).Value = 'Hortaliza'

# Row 514
$ws.Cells.Item(514, 1).Value = 7
$ws.Cells.Item(514, 2).Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(514, 3).Value = 'Ñuble'
$ws.Cells.Item(514, 4).Value = Get-Date -Year 2021 -Month 4 -Day 20 -Hour 0 -Minute 0 -Second 0
$ws.Cells.Item(514, 5).Value = 16
$ws.Cells.Item(514, 6).Value = 100114014
$ws.Cells.Item(514, 7).Value = 'Betarraga'
$ws.Cells.Item(514, 8).Value = 'Sin especificar'
$ws.Cells.Item(514, 9).Value = 'Primera'
$ws.Cells.Item(514, 10).Value = 160
$ws.Cells.Item(514, 11).Value = 600
$ws.Cells.Item(514, 12).Value = 650
$ws.Cells.Item(514, 13).Value = 625
$ws.Cells.Item(514, 14).Value = '$/paquete 5 unidades'
$ws.Cells.Item(514, 15).Value = 'Provincia de Diguillín'
$ws.Cells.Item(514, 16).Value = 125
$ws.Cells.Item(514, 17).Value = 5
$ws.Cells.Item(514, 18).Value = 'Hortaliza'

# Row 515
$ws.Cells.Item(515, 1).Value = 7
$ws.Cells.Item(515, 2).Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(515, 3).Value = 'Ñuble'
$ws.Cells.Item(515, 4).Value = Get-Date -Year 2021 -Month 6 -Day 9 -Hour 0 -Minute 0 -Second 0
$ws.Cells.Item(515, 5).Value = 16
$ws.Cells.Item(515, 6).Value = 100114014
$ws.Cells.Item(515, 7).Value = 'Betarraga'
$ws.Cells.Item(515, 8).Value = 'Sin especificar'
$ws.Cells.Item(515, 9).Value = 'Primera'
$ws.Cells.Item(515, 10).Value = 120
$ws.Cells.Item(515, 11).Value = 600
$ws.Cells.Item(515, 12).Value = 650
$ws.Cells.Item(515, 13).Value = 625
$ws.Cells.Item(515, 14).Value = '$/paquete 5 unidades'
$ws.Cells.Item(515, 15).Value = 'Provincia de Diguillín'
$ws.Cells.Item(515, 16).Value = 125
$ws.Cells.Item(515, 17).Value = 5
$ws.Cells.Item(515, 18).Value = 'Hortaliza'

# Row 516
$ws.Cells.Item(516, 1).Value = 7
$ws.Cells.Item(516, 2).Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(516, 3).Value = 'Ñuble'
$ws.Cells.Item(516, 4).Value = Get-Date -Year 2021 -Month 6 -Day 9 -Hour 0 -Minute 0 -Second 0
$ws.Cells.Item(516, 5).Value = 16
$ws.Cells.Item(516, 6).Value = 100114014
$ws.Cells.Item(516, 7).Value = 'Betarraga'
$ws.Cells.Item(516, 8).Value = 'Sin especificar'
$ws.Cells.Item(516, 9).Value = 'Segunda'
$ws.Cells.Item(516, 10).Value = 120
$ws.Cells.Item(516, 11).Value = 500
$ws.Cells.Item(516, 12).Value = 550
$ws.Cells.Item(516, 13).Value = 525
$ws.Cells.Item(516, 14).Value = '$/paquete 5 unidades'
$ws.Cells.Item(516, 15).Value = 'Provincia de Diguillín'
$ws.Cells.Item(516, 16).Value = 105
$ws.Cells.Item(516, 17).Value = 5
$ws.Cells.Item(516, 18).Value = 'Hortaliza'

# Row 517
$ws.Cells.Item(517, 1).Value = 7
$ws.Cells.Item(517, 2).Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(517, 3).Value = 'Ñuble'
$ws.Cells.Item(517, 4).Value = Get-Date -Year 2022 -Month 9 -Day 8 -Hour 0 -Minute 0 -Second 0
$ws.Cells.Item(517, 5).Value = 16
$ws.Cells.Item(517, 6).Value = 100114014
$ws.Cells.Item(517, 7).Value = 'Betarraga'
$ws.Cells.Item(517, 8).Value = 'Sin especificar'
$ws.Cells.Item(517, 9).Value = 'Segunda'
$ws.Cells.Item(517, 10).Value = 150
$ws.Cells.Item(517, 11).Value = 800
$ws.Cells.Item(517, 12).Value = 800
$ws.Cells.Item(517, 13).Value = 800
$ws.Cells.Item(517, 14).Value = '$/paquete 5 unidades'
$ws.Cells.Item(517, 15).Value = 'Provincia de Diguillín'
$ws.Cells.Item(517, 16).Value = 160
$ws.Cells.Item(517, 17).Value = 5
$ws.Cells.Item(517, 18).Value = 'Hortaliza'

# Row 518
$ws.Cells.Item(518, 1).Value = 7
$ws.Cells.Item(518, 2).Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(518, 3).Value = 'Ñuble'
$ws.Cells.Item(518, 4).Value = Get-Date -Year 2021 -Month 5 -Day 19 -Hour 0 -Minute 0 -Second 0
$ws.Cells.Item(518, 5).Value = 16
$ws.Cells.Item(518, 6).Value = 100114014
$ws.Cells.Item(518, 7).Value = 'Betarraga'
$ws.Cells.Item(518, 8).Value = 'Sin especificar'
$ws.Cells.Item(518, 9).Value = 'Primera'
$ws.Cells.Item(518, 10).Value = 300
$ws.Cells.Item(518, 11).Value = 600
$ws.Cells.Item(518, 12).Value = 650
$ws.Cells.Item(518, 13).Value = 625
$ws.Cells.Item(518, 14).Value = '$/paquete 5 unidades'
$ws.Cells.Item(518, 15).Value = 'Provincia de Diguillín'
$ws.Cells.Item(518, 16).Value = 125
$ws.Cells.Item(518, 17).Value = 5
$ws.Cells.Item(518, 18).Value = 'Hortaliza'

# Row 519
$ws.Cells.Item(519, 1).Value = 7
$ws.Cells.Item(519, 2).Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(519, 3).Value = 'Ñuble'
$ws.Cells.Item(519, 4).Value = Get-Date -Year 2021 -Month 5 -Day 19 -Hour 0 -Minute 0 -Second 0
$ws.Cells.Item(519, 5).Value = 16
$ws.Cells.Item(519, 6).Value = 100114014
$ws.Cells.Item(519, 7).Value = 'Betarraga'
$ws.Cells.Item(519, 8).Value = 'Sin especificar'
$ws.Cells.Item(519, 9).Value = 'Segunda'
$ws.Cells.Item(519, 10).Value = 120
$ws.Cells.Item(519, 11).Value = 500
$ws.Cells.Item(519, 12).Value = 550
$ws.Cells.Item(519, 13).Value = 525
$ws.Cells.Item(519, 14).Value = '$/paquete 5 unidades'
$ws.Cells.Item(519, 15).Value = 'Provincia de Diguillín'
$ws.Cells.Item(519, 16).Value = 105
$ws.Cells.Item(519, 17).Value = 5
$ws.Cells.Item(519, 18).Value = 'Hortaliza'

